# Backup QR Scanner data - 14/08/2025, 8:29:26 AM
#
# The scanner log previously held two rows (IDs 555586 and 452255). The
# older entry (555586, logged 08:28:49) is removed and the newer entry
# (452255, logged 08:28:54) becomes row 2. The worksheet is also renamed
# from "Scanner" to "Anatomy" to reflect the subject being logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the subject.
$ws.Name = "Anatomy"

# Drop the oldest scan (row 2, student 555586 @ 08:28:49). Deleting the
# row shifts the remaining row (previously row 3, student 452255 @
# 08:28:54) up into row 2, and the sheet's used range shrinks from
# A1:F3 to A1:F2 automatically.
$ws.Rows(2).Delete()
